$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add column C (13-01-2023), matching B1 style/border/bold.
# Use Value2 (not Value) so the date-like text is not auto-converted to a date serial.
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("C1").Value2 = "13-01-2023"

# Fund rows, "avg" and "total" reordered to the bottom, with new column C values added.
$ws.Range("A2").Value2 = "1822 Raices Valores Negociables"
$ws.Range("B2").Value2 = 50366.13
$ws.Range("C2").Value2 = 50309.35

$ws.Range("A3").Value2 = "Adcap IOL Acciones Argentina"
$ws.Range("B3").Value2 = 2379.69
$ws.Range("C3").Value2 = 2383.3

$ws.Range("A4").Value2 = "Allaria Acciones"
$ws.Range("B4").Value2 = 6511.13
$ws.Range("C4").Value2 = 6537.3

$ws.Range("A5").Value2 = "Alpha Mega"
$ws.Range("B5").Value2 = 18998.69
$ws.Range("C5").Value2 = 19092.07

$ws.Range("A6").Value2 = "Alpha renta balan global"
$ws.Range("B6").Value2 = 8085.95
$ws.Range("C6").Value2 = 9069.309999999999

$ws.Range("A7").Value2 = "Balanz"
$ws.Range("B7").Value2 = 18483.8
$ws.Range("C7").Value2 = 18544.26

$ws.Range("A8").Value2 = "Consultatio Acciones Argentina"
$ws.Range("B8").Value2 = 68831.62
$ws.Range("C8").Value2 = 68919.09

$ws.Range("A9").Value2 = "Consultatio Renta Variable"
$ws.Range("B9").Value2 = 4368.28
$ws.Range("C9").Value2 = 4522.63

$ws.Range("A10").Value2 = "Delta Acciones"
$ws.Range("B10").Value2 = 22403.68
$ws.Range("C10").Value2 = 22441.18

$ws.Range("A11").Value2 = "Delta Internacional"
$ws.Range("B11").Value2 = 1499.41
$ws.Range("C11").Value2 = 1500.71

$ws.Range("A12").Value2 = "Delta Latinoamerica"
$ws.Range("B12").Value2 = 1999.3
$ws.Range("C12").Value2 = 2002.88

$ws.Range("A13").Value2 = "Delta Select"
$ws.Range("B13").Value2 = 165125.55
$ws.Range("C13").Value2 = 164767.13

$ws.Range("A14").Value2 = "Delta gestion V"
$ws.Range("B14").Value2 = 37898.77
$ws.Range("C14").Value2 = 37960.32

$ws.Range("A15").Value2 = "FBA Acciones Argentinas"
$ws.Range("B15").Value2 = 30554.19
$ws.Range("C15").Value2 = 33459.83

$ws.Range("A16").Value2 = "FBA Calificado"
$ws.Range("B16").Value2 = 30759.08
$ws.Range("C16").Value2 = 31794.13

$ws.Range("A17").Value2 = "Fima Acciones"
$ws.Range("B17").Value2 = 40873.12
$ws.Range("C17").Value2 = 40649.77

$ws.Range("A18").Value2 = "Fima PB Acciones"
$ws.Range("B18").Value2 = 33306.76
$ws.Range("C18").Value2 = 33383.37

$ws.Range("A19").Value2 = "Gainvest Renta Variable"
$ws.Range("B19").Value2 = 55950.03
$ws.Range("C19").Value2 = 56011.95

$ws.Range("A20").Value2 = "Goal Acciones Argentinas"
$ws.Range("B20").Value2 = 2695.08
$ws.Range("C20").Value2 = 2677.9

$ws.Range("A21").Value2 = "Goal acciones plus"
$ws.Range("B21").Value2 = 403.28
$ws.Range("C21").Value2 = 404.84

$ws.Range("A22").Value2 = "HF Acciones Argentinas"
$ws.Range("B22").Value2 = 3706.03
$ws.Range("C22").Value2 = 3689.54

$ws.Range("A23").Value2 = "HF Acciones Lideres"
$ws.Range("B23").Value2 = 23177.7
$ws.Range("C23").Value2 = 23156.3

$ws.Range("A24").Value2 = "IEB Value"
$ws.Range("B24").Value2 = 1258.44
$ws.Range("C24").Value2 = 1258.12

$ws.Range("A25").Value2 = "Lombardi"
$ws.Range("B25").Value2 = 5990.75
$ws.Range("C25").Value2 = 6017.62

$ws.Range("A26").Value2 = "MAF"
$ws.Range("B26").Value2 = 5649.08
$ws.Range("C26").Value2 = 5627.43

$ws.Range("A27").Value2 = "Pellegrini Acciones"
$ws.Range("B27").Value2 = 18238.66
$ws.Range("C27").Value2 = 18282.58

$ws.Range("A28").Value2 = "Pionero Acciones"
$ws.Range("B28").Value2 = 7691.68
$ws.Range("C28").Value2 = 7602.22

$ws.Range("A29").Value2 = "Premier Renta Variable"
$ws.Range("B29").Value2 = 8934.290000000001
$ws.Range("C29").Value2 = 8916.5

$ws.Range("A30").Value2 = "Supefondo RV"
$ws.Range("B30").Value2 = 353937.75
$ws.Range("C30").Value2 = 353831.64

$ws.Range("A31").Value2 = "Superfondo "
$ws.Range("B31").Value2 = 387050.49
$ws.Range("C31").Value2 = 387147.47

$ws.Range("A32").Value2 = "Supergestion"
$ws.Range("B32").Value2 = 75667.28
$ws.Range("C32").Value2 = 75569.42

$ws.Range("A33").Value2 = "Toronto Trust Multimercado"
$ws.Range("B33").Value2 = 2696.38
$ws.Range("C33").Value2 = 2705.63

$ws.Range("A34").Value2 = "avg"
$ws.Range("B34").Value2 = 46734.13
$ws.Range("C34").Value2 = 46882.37

$ws.Range("A35").Value2 = "total"
$ws.Range("B35").Value2 = 1495492.07
$ws.Range("C35").Value2 = 1500235.79
